$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert a brand new sheet named "2022-Q1" right before "总计"
# (总计 currently sits at position 6; Worksheets.Add(Before) inserts before it)
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(6)
$newSheet = $wb.Worksheets.Add($zongji)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other quarterly sheets (0.75/0.75/1/1/.5/.5 in).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Step 2: seed "2022-Q1" with the same header/row formatting as "2021-Q4"
# (position 5 is unaffected by the insert above, since it sits before it)
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(5)
$q4.Range("A1:H14").Copy($newSheet.Range("A1"))

# Force the numeric-looking text columns (B,D,E,F,G) to stay text, matching
# the source data's storage (fund codes / percentages stored as strings).
$newSheet.Range("B2:B14").NumberFormat = "@"
$newSheet.Range("D2:F14").NumberFormat = "@"

# Row-by-row 2022-Q1 holdings data.
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "519017"
$newSheet.Cells.Item(2,3).Value = "大成积极成长混合"
$newSheet.Cells.Item(2,4).Value = "17.32"
$newSheet.Cells.Item(2,5).Value = "87.40"
$newSheet.Cells.Item(2,6).Value = "3.89"
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value = "0.6737"
$newSheet.Cells.Item(2,8).Value = 10

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "010490"
$newSheet.Cells.Item(3,3).Value = "鹏华高质量增长混合A"
$newSheet.Cells.Item(3,4).Value = "13.31"
$newSheet.Cells.Item(3,5).Value = "93.61"
$newSheet.Cells.Item(3,6).Value = "2.75"
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value = "0.3660"
$newSheet.Cells.Item(3,8).Value = 10

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "580002"
$newSheet.Cells.Item(4,3).Value = "东吴价值成长双动力混合A"
$newSheet.Cells.Item(4,4).Value = "2.89"
$newSheet.Cells.Item(4,5).Value = "85.87"
$newSheet.Cells.Item(4,6).Value = "2.87"
$newSheet.Cells.Item(4,7).NumberFormat = "@"
$newSheet.Cells.Item(4,7).Value = "0.0829"
$newSheet.Cells.Item(4,8).Value = 8

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "001415"
$newSheet.Cells.Item(5,3).Value = "信诚新锐回报灵活配置混合A"
$newSheet.Cells.Item(5,4).Value = "9.07"
$newSheet.Cells.Item(5,5).Value = "24.72"
$newSheet.Cells.Item(5,6).Value = "0.73"
$newSheet.Cells.Item(5,7).NumberFormat = "@"
$newSheet.Cells.Item(5,7).Value = "0.0662"
$newSheet.Cells.Item(5,8).Value = 5

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "001402"
$newSheet.Cells.Item(6,3).Value = "信诚新选回报灵活配置混合A"
$newSheet.Cells.Item(6,4).Value = "8.37"
$newSheet.Cells.Item(6,5).Value = "22.05"
$newSheet.Cells.Item(6,6).Value = "0.64"
$newSheet.Cells.Item(6,7).NumberFormat = "@"
$newSheet.Cells.Item(6,7).Value = "0.0536"
$newSheet.Cells.Item(6,8).Value = 4

$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "003234"
$newSheet.Cells.Item(7,3).Value = "信诚至利灵活配置混合A"
$newSheet.Cells.Item(7,4).Value = "8.99"
$newSheet.Cells.Item(7,5).Value = "22.05"
$newSheet.Cells.Item(7,6).Value = "0.59"
$newSheet.Cells.Item(7,7).NumberFormat = "@"
$newSheet.Cells.Item(7,7).Value = "0.0530"
$newSheet.Cells.Item(7,8).Value = 7

$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "004157"
$newSheet.Cells.Item(8,3).Value = "信诚至诚灵活配置混合A"
$newSheet.Cells.Item(8,4).Value = "7.32"
$newSheet.Cells.Item(8,5).Value = "22.71"
$newSheet.Cells.Item(8,6).Value = "0.68"
$newSheet.Cells.Item(8,7).NumberFormat = "@"
$newSheet.Cells.Item(8,7).Value = "0.0498"
$newSheet.Cells.Item(8,8).Value = 4

$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "002046"
$newSheet.Cells.Item(9,3).Value = "信诚新锐回报灵活配置混合B"
$newSheet.Cells.Item(9,4).Value = "5.19"
$newSheet.Cells.Item(9,5).Value = "24.72"
$newSheet.Cells.Item(9,6).Value = "0.73"
$newSheet.Cells.Item(9,7).NumberFormat = "@"
$newSheet.Cells.Item(9,7).Value = "0.0379"
$newSheet.Cells.Item(9,8).Value = 5

$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "003235"
$newSheet.Cells.Item(10,3).Value = "信诚至利灵活配置混合C"
$newSheet.Cells.Item(10,4).Value = "5.30"
$newSheet.Cells.Item(10,5).Value = "22.05"
$newSheet.Cells.Item(10,6).Value = "0.59"
$newSheet.Cells.Item(10,7).NumberFormat = "@"
$newSheet.Cells.Item(10,7).Value = "0.0313"
$newSheet.Cells.Item(10,8).Value = 7

$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "002030"
$newSheet.Cells.Item(11,3).Value = "信诚新选回报灵活配置混合B"
$newSheet.Cells.Item(11,4).Value = "3.85"
$newSheet.Cells.Item(11,5).Value = "22.05"
$newSheet.Cells.Item(11,6).Value = "0.64"
$newSheet.Cells.Item(11,7).NumberFormat = "@"
$newSheet.Cells.Item(11,7).Value = "0.0246"
$newSheet.Cells.Item(11,8).Value = 4

$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "004158"
$newSheet.Cells.Item(12,3).Value = "信诚至诚灵活配置混合B"
$newSheet.Cells.Item(12,4).Value = "2.18"
$newSheet.Cells.Item(12,5).Value = "22.71"
$newSheet.Cells.Item(12,6).Value = "0.68"
$newSheet.Cells.Item(12,7).NumberFormat = "@"
$newSheet.Cells.Item(12,7).Value = "0.0148"
$newSheet.Cells.Item(12,8).Value = 4

$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "010491"
$newSheet.Cells.Item(13,3).Value = "鹏华高质量增长混合C"
$newSheet.Cells.Item(13,4).Value = "0.28"
$newSheet.Cells.Item(13,5).Value = "93.61"
$newSheet.Cells.Item(13,6).Value = "2.75"
$newSheet.Cells.Item(13,7).NumberFormat = "@"
$newSheet.Cells.Item(13,7).Value = "0.0077"
$newSheet.Cells.Item(13,8).Value = 10

$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "011241"
$newSheet.Cells.Item(14,3).Value = "东吴价值成长双动力混合C"
$newSheet.Cells.Item(14,4).Value = "0.00"
$newSheet.Cells.Item(14,5).Value = "85.87"
$newSheet.Cells.Item(14,6).Value = "2.87"
$newSheet.Cells.Item(14,7).Value = 0
$newSheet.Cells.Item(14,8).Value = 8

# ---------------------------------------------------------------------------
# Step 3: update "总计" (now shifted to position 7) with the new 2022-Q1 row.
# Shift existing data rows (2-6) down to (3-7), bottom-up to avoid clobbering,
# then fill row 2 with the 2022-Q1 summary figures.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(7)
$total.Range("A6:D6").Copy($total.Range("A7"))
$total.Range("A5:D5").Copy($total.Range("A6"))
$total.Range("A4:D4").Copy($total.Range("A5"))
$total.Range("A3:D3").Copy($total.Range("A4"))
$total.Range("A2:D2").Copy($total.Range("A3"))

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 13
$total.Cells.Item(2,4).Value = 1.46

# Re-sequence the index column (A) for the rows that shifted down.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5

# Restore the original active sheet/tab (2020-Q4 was selected before this edit).
$wb.Worksheets.Item(1).Activate()
$newSheet.Range("A1").Select()
